$wb = $excel.ActiveWorkbook

# --- helper: write a value into a cell as TEXT (inline/shared string), keeping
#     the cell's style at the default (index 0) -- mirrors how the source
#     sheets store these "numeric-looking" columns as text.
$blankFormatSrc = $null
function Set-TextValue($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $script:blankFormatSrc.Copy()
    $cell.PasteSpecial(-4122)
}

# --- helper: write a real number into a cell, default style.
function Set-NumberValue($ws, $row, $col, $num) {
    $ws.Cells.Item($row, $col).Value = $num
}

$fundRows = @(
    @{A=0; B='040015'; C='华安动态灵活配置混合'; D='22.93'; E='79.55'; F='5.00'; G='1.1465'; H=2}
    @{A=1; B='010792'; C='华安成长先锋混合A'; D='16.44'; E='93.43'; F='6.32'; G='1.0390'; H=3}
    @{A=2; B='005299'; C='万家成长优选灵活配置混合A'; D='24.43'; E='93.91'; F='3.03'; G='0.7402'; H=7}
    @{A=3; B='006154'; C='华安制造先锋混合'; D='14.05'; E='93.81'; F='5.16'; G='0.7250'; H=3}
    @{A=4; B='010611'; C='万家战略发展产业混合A'; D='12.80'; E='93.73'; F='4.18'; G='0.5350'; H=6}
    @{A=5; B='010694'; C='万家内需增长一年持有期混合'; D='17.21'; E='94.85'; F='3.04'; G='0.5232'; H=10}
    @{A=6; B='501075'; C='万家科创主题3年封闭运作灵活配置混合A'; D='17.95'; E='79.41'; F='2.78'; G='0.4990'; H=8}
    @{A=7; B='040001'; C='华安创新混合'; D='16.06'; E='74.45'; F='3.03'; G='0.4866'; H=4}
    @{A=8; B='005300'; C='万家成长优选灵活配置混合C'; D='10.93'; E='93.91'; F='3.03'; G='0.3312'; H=7}
    @{A=9; B='010793'; C='华安成长先锋混合C'; D='3.80'; E='93.43'; F='6.32'; G='0.2402'; H=3}
    @{A=10; B='010612'; C='万家战略发展产业混合C'; D='5.06'; E='93.73'; F='4.18'; G='0.2115'; H=6}
    @{A=11; B='002707'; C='摩根士丹利华鑫科技领先灵活配置混合'; D='2.27'; E='93.05'; F='6.66'; G='0.1512'; H=3}
    @{A=12; B='501219'; C='华夏智胜先锋股票（LOF）A'; D='3.61'; E='94.50'; F='1.06'; G='0.0383'; H=7}
    @{A=13; B='006165'; C='建信中证1000指数增强A'; D='2.75'; E='93.00'; F='0.93'; G='0.0256'; H=8}
    @{A=14; B='014198'; C='华夏智胜先锋股票（LOF）C'; D='1.30'; E='94.50'; F='1.06'; G='0.0138'; H=7}
    @{A=15; B='080007'; C='长盛同鑫行业配置混合'; D='0.27'; E='87.31'; F='3.34'; G='0.0090'; H=8}
    @{A=16; B='006166'; C='建信中证1000指数增强C'; D='0.65'; E='93.00'; F='0.93'; G='0.0060'; H=8}
    @{A=17; B='013442'; C='建信中证1000指数增强E'; D='0.02'; E='93.00'; F='0.93'; G='0.0002'; H=8}
    @{A=18; B='007501'; C='万家科创主题3年封闭运作灵活配置混合C'; D=$null; E='79.41'; F='2.78'; G=0; H=8}
)
$totalRows = @(
    @{A=0; B='2022-Q1'; C=19; D=6.72}
    @{A=1; B='2021-Q4'; C=39; D=10.75}
    @{A=2; B='2021-Q3'; C=12; D=5.98}
    @{A=3; B='2021-Q2'; C=9; D=3.41}
    @{A=4; B='2021-Q1'; C=2; D=0.73}
)
# 1. Drop the old "总计" sheet -- it gets fully rebuilt (new top row + all
#    existing rows shift down), so it's simplest to recreate it from scratch.
$q4 = $wb.Worksheets.Item("2021-Q4")
$oldTotal = $wb.Worksheets.Item("总计")
$oldTotal.Delete()

# 2. Insert the new "2022-Q1" sheet right after "2021-Q4" (where "总计" used
#    to be), so sheet order stays 2021-Q1..Q4, 2022-Q1, 总计.
$q1sheet = $wb.Worksheets.Add($null, $q4)
$q1sheet.Name = "2022-Q1"

$blankFormatSrc = $q1sheet.Range("Z1")
$script:blankFormatSrc = $blankFormatSrc

# Reuse the existing header / index-column formatting (bold, centered,
# thin-bordered "style 2") from the 2021-Q4 sheet instead of re-building it
# from scratch, so the same style slot gets reused rather than a new one
# minted.
$q4.Range("B1:H1").Copy()
$q1sheet.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2").Copy()
$q1sheet.Range("A2:A20").PasteSpecial(-4122)

$q1sheet.Cells.Item(1,2).Value = "基金代码"
$q1sheet.Cells.Item(1,3).Value = "基金名称"
$q1sheet.Cells.Item(1,4).Value = "基金规模"
$q1sheet.Cells.Item(1,5).Value = "股票总仓位"
$q1sheet.Cells.Item(1,6).Value = "仓位占比"
$q1sheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1sheet.Cells.Item(1,8).Value = "仓位排名"

foreach ($r in $fundRows) {
    $row = $r.A + 2
    Set-NumberValue $q1sheet $row 1 $r.A
    Set-TextValue   $q1sheet $row 2 $r.B
    Set-TextValue   $q1sheet $row 3 $r.C
    if ($null -ne $r.D) {
        Set-TextValue $q1sheet $row 4 $r.D
    }
    Set-TextValue   $q1sheet $row 5 $r.E
    Set-TextValue   $q1sheet $row 6 $r.F
    if ($r.G -is [string]) {
        Set-TextValue $q1sheet $row 7 $r.G
    } else {
        Set-NumberValue $q1sheet $row 7 $r.G
    }
    Set-NumberValue $q1sheet $row 8 $r.H
}

# 3. Re-create "总计" after "2022-Q1", with the new quarter's row at the top.
$totalSheet = $wb.Worksheets.Add($null, $q1sheet)
$totalSheet.Name = "总计"

$q4.Range("B1:D1").Copy()
$totalSheet.Range("B1:D1").PasteSpecial(-4122)
$q4.Range("A2").Copy()
$totalSheet.Range("A2:A6").PasteSpecial(-4122)

$totalSheet.Cells.Item(1,2).Value = "日期"
$totalSheet.Cells.Item(1,3).Value = "持有数量(只)"
$totalSheet.Cells.Item(1,4).Value = "持有市值(亿元)"

foreach ($r in $totalRows) {
    $row = $r.A + 2
    Set-NumberValue $totalSheet $row 1 $r.A
    Set-TextValue   $totalSheet $row 2 $r.B
    Set-NumberValue $totalSheet $row 3 $r.C
    Set-NumberValue $totalSheet $row 4 $r.D
}

Write-Output "Sheets now: $($wb.Worksheets.Count)"
